# Changed Card Numbers for FeeData
# ---------------------------------------------------------------------------
# Adds a new "SEC" (SEC Code) column to Sheet3's lookup table and stamps a
# "Public" classification footer on every worksheet.

$wb = $excel.ActiveWorkbook

# --- Sheet3: add new column AA ("SEC") with per-row SEC code values -------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()

$ws3.Range("AA1").Value = "SEC"
$ws3.Range("AA2").Value = "TEL"
$ws3.Range("AA3").Value = "TEL"
$ws3.Range("AA4").Value = "CCD"

# Move the selection to reflect the newly-entered data (matches the
# workbook's last on-screen selection after the edit).
$ws3.Range("AA4").Select()

# --- Classification footer on every worksheet ------------------------------
$footerText = [char]13 + "&1#&`"Calibri`"&10&K000000 Public "

foreach ($sheetName in @("Sheet1", "Sheet2", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ps = $ws.PageSetup
    $ps.CenterFooter = $footerText
}
